# Auto-generated Excel COM-interop script applying scheduled-runner price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 1994.5946
$ws.Cells.Item(17, 10).Value = 2183.3333
$ws.Cells.Item(17, 12).Value = 6549.999899999999
$ws.Cells.Item(17, 14).Value = -6885.999899999999
$ws.Cells.Item(20, 8).Value = 2000
$ws.Cells.Item(20, 9).Value = 2000
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 2000
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = -1770
$ws.Cells.Item(20, 14).ClearContents()
$ws.Cells.Item(35, 8).Value = 2000
$ws.Cells.Item(35, 9).Value = 2000
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 2000
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -1621
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(55, 8).Value = 666
$ws.Cells.Item(55, 9).Value = 430.54544
$ws.Cells.Item(55, 10).Value = 1036
$ws.Cells.Item(55, 11).Value = 430.54544
$ws.Cells.Item(55, 12).Value = 1036
$ws.Cells.Item(55, 13).Value = -216.54544
$ws.Cells.Item(55, 14).Value = -1464
$ws.Cells.Item(100, 8).Value = 4005
$ws.Cells.Item(100, 10).Value = 4306.25
$ws.Cells.Item(100, 12).Value = 4306.25
$ws.Cells.Item(100, 14).Value = -5388.25
$ws.Cells.Item(116, 8).Value = 9316.083000000001
$ws.Cells.Item(116, 9).Value = 8599
$ws.Cells.Item(116, 11).Value = 8599
$ws.Cells.Item(116, 13).Value = -5157
$ws.Cells.Item(141, 8).Value = 5884.4614
$ws.Cells.Item(141, 9).Value = 874.5
$ws.Cells.Item(141, 10).Value = 8111.1113
$ws.Cells.Item(141, 11).Value = 2623.5
$ws.Cells.Item(141, 12).Value = 24333.3339
$ws.Cells.Item(141, 13).Value = 2556.5
$ws.Cells.Item(141, 14).Value = -34693.3339

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 8).Value = 249.91667
$ws.Cells.Item(5, 9).Value = 167.66667
$ws.Cells.Item(5, 10).Value = 496.66666
$ws.Cells.Item(5, 11).Value = 167.66667
$ws.Cells.Item(5, 12).Value = 496.66666
$ws.Cells.Item(5, 13).Value = -55.66667000000001
$ws.Cells.Item(5, 14).Value = -720.66666
$ws.Cells.Item(26, 8).Value = 385.66666
$ws.Cells.Item(26, 9).Value = 353.5
$ws.Cells.Item(26, 10).Value = 450
$ws.Cells.Item(26, 11).Value = 353.5
$ws.Cells.Item(26, 12).Value = 450
$ws.Cells.Item(26, 13).Value = -23.5
$ws.Cells.Item(26, 14).Value = -1110
$ws.Cells.Item(39, 8).Value = 1500
$ws.Cells.Item(39, 10).Value = 2000
$ws.Cells.Item(39, 12).Value = 2000
$ws.Cells.Item(39, 14).Value = -3040
$ws.Cells.Item(45, 8).Value = 111116240
$ws.Cells.Item(45, 9).Value = 166667860
$ws.Cells.Item(45, 10).Value = 13004.333
$ws.Cells.Item(45, 11).Value = 166667860
$ws.Cells.Item(45, 12).Value = 13004.333
$ws.Cells.Item(45, 13).Value = -166667483
$ws.Cells.Item(45, 14).Value = -13758.333
$ws.Cells.Item(74, 8).Value = 25644120
$ws.Cells.Item(74, 9).Value = 41670450
$ws.Cells.Item(74, 11).Value = 41670450
$ws.Cells.Item(74, 13).Value = -41669576
$ws.Cells.Item(77, 8).Value = 25644120
$ws.Cells.Item(77, 9).Value = 41670450
$ws.Cells.Item(77, 11).Value = 208352250
$ws.Cells.Item(77, 13).Value = -208347882
$ws.Cells.Item(88, 8).Value = 3465.4666
$ws.Cells.Item(88, 10).Value = 2952.9092
$ws.Cells.Item(88, 12).Value = 2952.9092
$ws.Cells.Item(88, 14).Value = -3764.9092
$ws.Cells.Item(91, 8).Value = 3465.4666
$ws.Cells.Item(91, 10).Value = 2952.9092
$ws.Cells.Item(91, 12).Value = 2952.9092
$ws.Cells.Item(91, 14).Value = -5760.9092
$ws.Cells.Item(132, 8).Value = 3257.6
$ws.Cells.Item(132, 9).Value = 2611.111
$ws.Cells.Item(132, 10).Value = 4920
$ws.Cells.Item(132, 11).Value = 7833.333
$ws.Cells.Item(132, 12).Value = 14760
$ws.Cells.Item(132, 13).Value = -5303.333
$ws.Cells.Item(132, 14).Value = -19820

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 8).Value = 249.91667
$ws.Cells.Item(4, 9).Value = 167.66667
$ws.Cells.Item(4, 10).Value = 496.66666
$ws.Cells.Item(4, 11).Value = 167.66667
$ws.Cells.Item(4, 12).Value = 496.66666
$ws.Cells.Item(4, 13).Value = -52.66667000000001
$ws.Cells.Item(4, 14).Value = -726.66666
$ws.Cells.Item(22, 8).Value = 227.33333
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(19, 8).Value = 2500260.5
$ws.Cells.Item(19, 9).Value = 5000029
$ws.Cells.Item(19, 10).Value = 492.25
$ws.Cells.Item(19, 11).Value = 5000029
$ws.Cells.Item(19, 12).Value = 492.25
$ws.Cells.Item(19, 13).Value = -4999859
$ws.Cells.Item(19, 14).Value = -832.25
$ws.Cells.Item(24, 8).Value = 2500260.5
$ws.Cells.Item(24, 9).Value = 5000029
$ws.Cells.Item(24, 10).Value = 492.25
$ws.Cells.Item(24, 11).Value = 5000029
$ws.Cells.Item(24, 12).Value = 492.25
$ws.Cells.Item(24, 13).Value = -4999859
$ws.Cells.Item(24, 14).Value = -832.25
$ws.Cells.Item(35, 8).Value = 8420
$ws.Cells.Item(35, 9).Value = 8189.2
$ws.Cells.Item(35, 10).Value = 10728
$ws.Cells.Item(35, 11).Value = 8189.2
$ws.Cells.Item(35, 12).Value = 10728
$ws.Cells.Item(35, 13).Value = -7895.2
$ws.Cells.Item(35, 14).Value = -11316

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 7645849
$ws.Cells.Item(4, 9).Value = 4375097.5
$ws.Cells.Item(4, 11).Value = 13125292.5
$ws.Cells.Item(4, 13).Value = -13125180.5
$ws.Cells.Item(35, 8).Value = 600
$ws.Cells.Item(35, 9).Value = 633.3333
$ws.Cells.Item(35, 10).Value = 583.3333
$ws.Cells.Item(35, 11).Value = 1899.9999
$ws.Cells.Item(35, 12).Value = 1749.9999
$ws.Cells.Item(35, 13).Value = -1611.9999
$ws.Cells.Item(35, 14).Value = -2325.9999
$ws.Cells.Item(52, 8).Value = 483
$ws.Cells.Item(52, 10).Value = 483
$ws.Cells.Item(52, 12).Value = 1449
$ws.Cells.Item(52, 14).Value = -1981
$ws.Cells.Item(86, 8).Value = 2117.818
$ws.Cells.Item(86, 10).Value = 3056.2856
$ws.Cells.Item(86, 12).Value = 9168.856800000001
$ws.Cells.Item(86, 14).Value = -11540.8568
$ws.Cells.Item(89, 8).Value = 2117.818
$ws.Cells.Item(89, 10).Value = 3056.2856
$ws.Cells.Item(89, 12).Value = 27506.5704
$ws.Cells.Item(89, 14).Value = -39362.5704

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 14).Value = -1224
$ws.Cells.Item(31, 8).Value = 3800
$ws.Cells.Item(31, 9).Value = 3800
$ws.Cells.Item(31, 11).Value = 3800
$ws.Cells.Item(31, 13).Value = -3508
$ws.Cells.Item(37, 8).Value = 3800
$ws.Cells.Item(37, 9).Value = 3800
$ws.Cells.Item(37, 11).Value = 3800
$ws.Cells.Item(37, 13).Value = -3523
$ws.Cells.Item(53, 8).Value = 13010.5
$ws.Cells.Item(53, 10).Value = 13010.5
$ws.Cells.Item(53, 12).Value = 13010.5
$ws.Cells.Item(53, 14).Value = -14272.5
$ws.Cells.Item(70, 8).Value = 17382.166
$ws.Cells.Item(70, 9).Value = 6375.6924
$ws.Cells.Item(70, 11).Value = 6375.6924
$ws.Cells.Item(70, 13).Value = -6105.6924
$ws.Cells.Item(73, 8).Value = 17382.166
$ws.Cells.Item(73, 9).Value = 6375.6924
$ws.Cells.Item(73, 11).Value = 6375.6924
$ws.Cells.Item(73, 13).Value = -5439.6924
$ws.Cells.Item(132, 8).Value = 73138.60000000001
$ws.Cells.Item(132, 9).Value = 128205.375
$ws.Cells.Item(132, 10).Value = 10205.143
$ws.Cells.Item(132, 11).Value = 384616.125
$ws.Cells.Item(132, 12).Value = 30615.429
$ws.Cells.Item(132, 13).Value = -382086.125
$ws.Cells.Item(132, 14).Value = -35675.429
$ws.Cells.Item(135, 8).Value = 69796.8
$ws.Cells.Item(135, 10).Value = 69796.8
$ws.Cells.Item(135, 12).Value = 69796.8
$ws.Cells.Item(135, 14).Value = -79936.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(9, 8).Value = 1130
$ws.Cells.Item(9, 9).Value = 1130
$ws.Cells.Item(9, 11).Value = 1130
$ws.Cells.Item(9, 13).Value = -906
$ws.Cells.Item(30, 8).Value = 2150.5
$ws.Cells.Item(30, 9).Value = 2300
$ws.Cells.Item(30, 10).Value = 2001
$ws.Cells.Item(30, 11).Value = 2300
$ws.Cells.Item(30, 12).Value = 2001
$ws.Cells.Item(30, 13).Value = -2192
$ws.Cells.Item(30, 14).Value = -2217
$ws.Cells.Item(132, 8).Value = 4279.857
$ws.Cells.Item(132, 9).Value = 2628.7896
$ws.Cells.Item(132, 10).Value = 7765.4443
$ws.Cells.Item(132, 11).Value = 7886.3688
$ws.Cells.Item(132, 12).Value = 23296.3329
$ws.Cells.Item(132, 13).Value = -5356.3688
$ws.Cells.Item(132, 14).Value = -28356.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 3445.2307
$ws.Cells.Item(81, 9).Value = 2225.75
$ws.Cells.Item(81, 10).Value = 5396.4
$ws.Cells.Item(81, 11).Value = 4451.5
$ws.Cells.Item(81, 12).Value = 10792.8
$ws.Cells.Item(81, 13).Value = -3390.5
$ws.Cells.Item(81, 14).Value = -12914.8
$ws.Cells.Item(84, 8).Value = 3445.2307
$ws.Cells.Item(84, 9).Value = 2225.75
$ws.Cells.Item(84, 10).Value = 5396.4
$ws.Cells.Item(84, 11).Value = 22257.5
$ws.Cells.Item(84, 12).Value = 53964
$ws.Cells.Item(84, 13).Value = -16953.5
$ws.Cells.Item(84, 14).Value = -64572
$ws.Cells.Item(98, 8).Value = 30000
$ws.Cells.Item(98, 10).Value = 30000
$ws.Cells.Item(98, 12).Value = 30000
$ws.Cells.Item(98, 14).Value = -35990
$ws.Cells.Item(110, 8).Value = 59995
$ws.Cells.Item(110, 10).Value = 59995
$ws.Cells.Item(110, 12).Value = 59995
$ws.Cells.Item(110, 14).Value = -68175
$ws.Cells.Item(132, 8).Value = 5403
$ws.Cells.Item(132, 9).Value = 4670.2856
$ws.Cells.Item(132, 10).Value = 9249.75
$ws.Cells.Item(132, 11).Value = 14010.8568
$ws.Cells.Item(132, 12).Value = 27749.25
$ws.Cells.Item(132, 13).Value = -11480.8568
$ws.Cells.Item(132, 14).Value = -32809.25
$ws.Cells.Item(133, 8).Value = 36277.332
$ws.Cells.Item(133, 10).Value = 36277.332
$ws.Cells.Item(133, 12).Value = 36277.332
$ws.Cells.Item(133, 14).Value = -46397.332
$ws.Cells.Item(136, 8).Value = 3754.4583
$ws.Cells.Item(136, 9).Value = 2054.85
$ws.Cells.Item(136, 10).Value = 12252.5
$ws.Cells.Item(136, 11).Value = 6164.549999999999
$ws.Cells.Item(136, 12).Value = 36757.5
$ws.Cells.Item(136, 13).Value = -3614.549999999999
$ws.Cells.Item(136, 14).Value = -41857.5
